$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-23 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-24 Monday", 2) | Out-Null
$d.Content.Find.Execute("76-73=3", $true, $false, $false, $false, $false, $true, 1, $false, "16+81=97", 2) | Out-Null
$d.Content.Find.Execute("56-54=2", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=38", 2) | Out-Null
$d.Content.Find.Execute("30+38=68", $true, $false, $false, $false, $false, $true, 1, $false, "60-54=6", 2) | Out-Null
$d.Content.Find.Execute("57-12=45", $true, $false, $false, $false, $false, $true, 1, $false, "21+28=49", 2) | Out-Null
$d.Content.Find.Execute("51+14=65", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=67", 2) | Out-Null
$d.Content.Find.Execute("3+23=26", $true, $false, $false, $false, $false, $true, 1, $false, "12+55=67", 2) | Out-Null
$d.Content.Find.Execute("38+11=49", $true, $false, $false, $false, $false, $true, 1, $false, "8+61=69", 2) | Out-Null
$d.Content.Find.Execute("9+3=12", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=59", 2) | Out-Null
$d.Content.Find.Execute("87-16=71", $true, $false, $false, $false, $false, $true, 1, $false, "7+12=19", 2) | Out-Null
$d.Content.Find.Execute("65+29=94", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=71", 2) | Out-Null
$d.Content.Find.Execute("95-34=61", $true, $false, $false, $false, $false, $true, 1, $false, "89-13=76", 2) | Out-Null
$d.Content.Find.Execute("16+15=31", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=29", 2) | Out-Null
$d.Content.Find.Execute("93-84=9", $true, $false, $false, $false, $false, $true, 1, $false, "99-48=51", 2) | Out-Null
$d.Content.Find.Execute("62-58=4", $true, $false, $false, $false, $false, $true, 1, $false, "53-0=53", 2) | Out-Null
$d.Content.Find.Execute("5+47=52", $true, $false, $false, $false, $false, $true, 1, $false, "35+29=64", 2) | Out-Null
$d.Content.Find.Execute("34-17=17", $true, $false, $false, $false, $false, $true, 1, $false, "37-3=34", 2) | Out-Null
$d.Content.Find.Execute("53-40=13", $true, $false, $false, $false, $false, $true, 1, $false, "17+82=99", 2) | Out-Null
$d.Content.Find.Execute("34+16=50", $true, $false, $false, $false, $false, $true, 1, $false, "19+45=64", 2) | Out-Null
$d.Content.Find.Execute("78-31=47", $true, $false, $false, $false, $false, $true, 1, $false, "84+11=95", 2) | Out-Null
$d.Content.Find.Execute("4+87=91", $true, $false, $false, $false, $false, $true, 1, $false, "64+16=80", 2) | Out-Null
$d.Content.Find.Execute("60-16=44", $true, $false, $false, $false, $false, $true, 1, $false, "25+53=78", 2) | Out-Null
$d.Content.Find.Execute("39+13=52", $true, $false, $false, $false, $false, $true, 1, $false, "65-50=15", 2) | Out-Null
$d.Content.Find.Execute("49-21=28", $true, $false, $false, $false, $false, $true, 1, $false, "65-33=32", 2) | Out-Null
$d.Content.Find.Execute("21+26=47", $true, $false, $false, $false, $false, $true, 1, $false, "28-24=4", 2) | Out-Null
$d.Content.Find.Execute("67-0=67", $true, $false, $false, $false, $false, $true, 1, $false, "93-71=22", 2) | Out-Null
$d.Content.Find.Execute("58-56=2", $true, $false, $false, $false, $false, $true, 1, $false, "35+51=86", 2) | Out-Null
$d.Content.Find.Execute("73-11=62", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=62", 2) | Out-Null
$d.Content.Find.Execute("76-5=71", $true, $false, $false, $false, $false, $true, 1, $false, "75-36=39", 2) | Out-Null
$d.Content.Find.Execute("81-80=1", $true, $false, $false, $false, $false, $true, 1, $false, "15-3=12", 2) | Out-Null
$d.Content.Find.Execute("97-4=93", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=7", 2) | Out-Null
$d.Content.Find.Execute("89-77=12", $true, $false, $false, $false, $false, $true, 1, $false, "69+13=82", 2) | Out-Null
$d.Content.Find.Execute("47-20=27", $true, $false, $false, $false, $false, $true, 1, $false, "44+26=70", 2) | Out-Null
$d.Content.Find.Execute("61-30=31", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=15", 2) | Out-Null
$d.Content.Find.Execute("25-11=14", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=81", 2) | Out-Null
$d.Content.Find.Execute("61-24=37", $true, $false, $false, $false, $false, $true, 1, $false, "14+62=76", 2) | Out-Null
$d.Content.Find.Execute("60-12=48", $true, $false, $false, $false, $false, $true, 1, $false, "85-41=44", 2) | Out-Null
$d.Content.Find.Execute("70+4=74", $true, $false, $false, $false, $false, $true, 1, $false, "22+47=69", 2) | Out-Null
$d.Content.Find.Execute("76+18=94", $true, $false, $false, $false, $false, $true, 1, $false, "42+54=96", 2) | Out-Null
$d.Content.Find.Execute("56+14=70", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=91", 2) | Out-Null
$d.Content.Find.Execute("68-22=46", $true, $false, $false, $false, $false, $true, 1, $false, "77+2=79", 2) | Out-Null
$d.Content.Find.Execute("8+55=63", $true, $false, $false, $false, $false, $true, 1, $false, "46+29=75", 2) | Out-Null
$d.Content.Find.Execute("77-15=62", $true, $false, $false, $false, $false, $true, 1, $false, "46-39=7", 2) | Out-Null
$d.Content.Find.Execute("60+36=96", $true, $false, $false, $false, $false, $true, 1, $false, "2+44=46", 2) | Out-Null
$d.Content.Find.Execute("61-13=48", $true, $false, $false, $false, $false, $true, 1, $false, "77-75=2", 2) | Out-Null
$d.Content.Find.Execute("82-74=8", $true, $false, $false, $false, $false, $true, 1, $false, "90-48=42", 2) | Out-Null
$d.Content.Find.Execute("10+8=18", $true, $false, $false, $false, $false, $true, 1, $false, "42+21=63", 2) | Out-Null
$d.Content.Find.Execute("31-12=19", $true, $false, $false, $false, $false, $true, 1, $false, "60-10=50", 2) | Out-Null
$d.Content.Find.Execute("66-43=23", $true, $false, $false, $false, $false, $true, 1, $false, "28+59=87", 2) | Out-Null
$d.Content.Find.Execute("63-30=33", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=38", 2) | Out-Null
$d.Content.Find.Execute("14+67=81", $true, $false, $false, $false, $false, $true, 1, $false, "11+60=71", 2) | Out-Null
$d.Content.Find.Execute("49-30=19", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=6", 2) | Out-Null
$d.Content.Find.Execute("28+66=94", $true, $false, $false, $false, $false, $true, 1, $false, "81-64=17", 2) | Out-Null
$d.Content.Find.Execute("51-16=35", $true, $false, $false, $false, $false, $true, 1, $false, "89-29=60", 2) | Out-Null
$d.Content.Find.Execute("97-95=2", $true, $false, $false, $false, $false, $true, 1, $false, "35-33=2", 2) | Out-Null
$d.Content.Find.Execute("89-71=18", $true, $false, $false, $false, $false, $true, 1, $false, "72-51=21", 2) | Out-Null
$d.Content.Find.Execute("64-61=3", $true, $false, $false, $false, $false, $true, 1, $false, "57-26=31", 2) | Out-Null
$d.Content.Find.Execute("7+74=81", $true, $false, $false, $false, $false, $true, 1, $false, "12+41=53", 2) | Out-Null
$d.Content.Find.Execute("90-6=84", $true, $false, $false, $false, $false, $true, 1, $false, "56-53=3", 2) | Out-Null
$d.Content.Find.Execute("76-9=67", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=55", 2) | Out-Null
$d.Content.Find.Execute("76+20=96", $true, $false, $false, $false, $false, $true, 1, $false, "21-9=12", 2) | Out-Null
$d.Content.Find.Execute("10-8=2", $true, $false, $false, $false, $false, $true, 1, $false, "80-41=39", 2) | Out-Null
$d.Content.Find.Execute("35-31=4", $true, $false, $false, $false, $false, $true, 1, $false, "93-82=11", 2) | Out-Null
$d.Content.Find.Execute("17+8=25", $true, $false, $false, $false, $false, $true, 1, $false, "72-5=67", 2) | Out-Null
$d.Content.Find.Execute("48-40=8", $true, $false, $false, $false, $false, $true, 1, $false, "66+21=87", 2) | Out-Null
$d.Content.Find.Execute("87-83=4", $true, $false, $false, $false, $false, $true, 1, $false, "32-23=9", 2) | Out-Null
$d.Content.Find.Execute("6+62=68", $true, $false, $false, $false, $false, $true, 1, $false, "81-25=56", 2) | Out-Null
$d.Content.Find.Execute("89-37=52", $true, $false, $false, $false, $false, $true, 1, $false, "9+8=17", 2) | Out-Null
$d.Content.Find.Execute("3+19=22", $true, $false, $false, $false, $false, $true, 1, $false, "99-25=74", 2) | Out-Null
$d.Content.Find.Execute("71-10=61", $true, $false, $false, $false, $false, $true, 1, $false, "70-18=52", 2) | Out-Null
$d.Content.Find.Execute("53-15=38", $true, $false, $false, $false, $false, $true, 1, $false, "88+2=90", 2) | Out-Null
$d.Content.Find.Execute("0-0=0", $true, $false, $false, $false, $false, $true, 1, $false, "11+28=39", 2) | Out-Null
$d.Content.Find.Execute("45+17=62", $true, $false, $false, $false, $false, $true, 1, $false, "98-61=37", 2) | Out-Null
$d.Content.Find.Execute("80+18=98", $true, $false, $false, $false, $false, $true, 1, $false, "77+10=87", 2) | Out-Null
$d.Content.Find.Execute("0+10=10", $true, $false, $false, $false, $false, $true, 1, $false, "19+22=41", 2) | Out-Null
$d.Content.Find.Execute("96-69=27", $true, $false, $false, $false, $false, $true, 1, $false, "90-88=2", 2) | Out-Null
$d.Content.Find.Execute("76-7=69", $true, $false, $false, $false, $false, $true, 1, $false, "22+4=26", 2) | Out-Null
$d.Content.Find.Execute("19+75=94", $true, $false, $false, $false, $false, $true, 1, $false, "65-60=5", 2) | Out-Null
$d.Content.Find.Execute("10+62=72", $true, $false, $false, $false, $false, $true, 1, $false, "85-66=19", 2) | Out-Null
$d.Content.Find.Execute("66-15=51", $true, $false, $false, $false, $false, $true, 1, $false, "71-41=30", 2) | Out-Null
$d.Content.Find.Execute("98-45=53", $true, $false, $false, $false, $false, $true, 1, $false, "74-4=70", 2) | Out-Null
$d.Content.Find.Execute("72+24=96", $true, $false, $false, $false, $false, $true, 1, $false, "3-2=1", 2) | Out-Null
$d.Content.Find.Execute("38+1=39", $true, $false, $false, $false, $false, $true, 1, $false, "6+84=90", 2) | Out-Null
$d.Content.Find.Execute("48+19=67", $true, $false, $false, $false, $false, $true, 1, $false, "63+26=89", 2) | Out-Null
$d.Content.Find.Execute("43-31=12", $true, $false, $false, $false, $false, $true, 1, $false, "28+11=39", 2) | Out-Null
$d.Content.Find.Execute("44+38=82", $true, $false, $false, $false, $false, $true, 1, $false, "27-19=8", 2) | Out-Null
$d.Content.Find.Execute("83+16=99", $true, $false, $false, $false, $false, $true, 1, $false, "12+5=17", 2) | Out-Null
$d.Content.Find.Execute("60-20=40", $true, $false, $false, $false, $false, $true, 1, $false, "27-11=16", 2) | Out-Null
$d.Content.Find.Execute("20+21=41", $true, $false, $false, $false, $false, $true, 1, $false, "2+77=79", 2) | Out-Null
$d.Content.Find.Execute("13+51=64", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=56", 2) | Out-Null
$d.Content.Find.Execute("9+7=16", $true, $false, $false, $false, $false, $true, 1, $false, "37+33=70", 2) | Out-Null
$d.Content.Find.Execute("86-4=82", $true, $false, $false, $false, $false, $true, 1, $false, "40-32=8", 2) | Out-Null
$d.Content.Find.Execute("73-37=36", $true, $false, $false, $false, $false, $true, 1, $false, "34+11=45", 2) | Out-Null
$d.Content.Find.Execute("0+79=79", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=39", 2) | Out-Null
$d.Content.Find.Execute("96-27=69", $true, $false, $false, $false, $false, $true, 1, $false, "52-19=33", 2) | Out-Null
$d.Content.Find.Execute("60-47=13", $true, $false, $false, $false, $false, $true, 1, $false, "73-9=64", 2) | Out-Null
$d.Content.Find.Execute("38+61=99", $true, $false, $false, $false, $false, $true, 1, $false, "9+16=25", 2) | Out-Null
$d.Content.Find.Execute("95-31=64", $true, $false, $false, $false, $false, $true, 1, $false, "4+10=14", 2) | Out-Null
$d.Content.Find.Execute("44-32=12", $true, $false, $false, $false, $false, $true, 1, $false, "46-23=23", 2) | Out-Null
$d.Content.Find.Execute("75-18=57", $true, $false, $false, $false, $false, $true, 1, $false, "95-94=1", 2) | Out-Null
$d.Content.Find.Execute("7+75=82", $true, $false, $false, $false, $false, $true, 1, $false, "47-42=5", 2) | Out-Null
